$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" values between row 16 and row 20 (F column)
$f16 = $ws.Range("F16").Value2
$f20 = $ws.Range("F20").Value2

$ws.Range("F16").Value2 = $f20
$ws.Range("F20").Value2 = $f16

$wb.Save()
